# Quest 8.docx edit script
# Applies:
#  1. Header restructure: "CPSC121 SI" / "By: Derek Louie" paragraph ->
#     "200 EXP" paragraph + new "CPSC121 SI" paragraph (carrying the
#     _GoBack bookmark, dropping the "By: Derek Louie" byline).
#  2-20. Clean up stray spell/grammar-check proofing marks that split
#     otherwise-contiguous sentences/code lines into multiple runs, by
#     doing a no-op Find & Replace over the full contiguous text (which
#     collapses the runs Word had split for the proofErr tags).
#  21. Final paragraph: merge runs around the relocated _GoBack bookmark
#     (bookmark no longer needed there since it now lives in the header).

$d = $word.ActiveDocument

$lq = [char]0x201C
$rq = [char]0x201D

function Replace-Literal($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
    return $ok
}

# 1. Header restructure --------------------------------------------------
Replace-Literal "CPSC121 SI" "200 EXP^pCPSC121 SI"
Replace-Literal ([char]11 + "By: Derek Louie") ""

$headerPara = $d.Paragraphs(3)
$bmStart = $d.Range($headerPara.Range.Start, $headerPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmStart)

# 2. Rex: nested if statement explanation --------------------------------
Replace-Literal "A nested if statement is an if statement inside of another if statement." "A nested if statement is an if statement inside of another if statement."

# 3. if(school == "CSUF") code line --------------------------------------
Replace-Literal ("    if(school == " + $lq + "CSUF" + $rq + ")") ("    if(school == " + $lq + "CSUF" + $rq + ")")

# 4. cout << "You are in CPSC121-07 at CSUF\n";
Replace-Literal ("cout << " + $lq + "You are in CPSC121-07 at CSUF\n" + $rq + ";") ("cout << " + $lq + "You are in CPSC121-07 at CSUF\n" + $rq + ";")

Write-Output "done so far"
